$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "42.896.24"
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.359.26"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +1.53%  "

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "302.80"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.34%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "95.15"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +1.26%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("E8").Value = "  -0.23%  "

$ws.Range("E9").Value = "  -3.49%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "34.25"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +1.02%  "

$ws.Range("E11").Value = "  +2.22%  "

$ws.Range("E12").Value = "  +0.48%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "18.40"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -1.91%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.71"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.25%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "2.725.80"
$c.Style = "Normal"

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "2.377.08"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +2.40%  "

$ws.Range("E17").Value = "  +0.77%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "42.913.44"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.44%  "

$ws.Range("E19").Value = "  +0.09%  "

$ws.Range("E20").Value = "  +1.81%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "0.0₃0884"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.16%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "67.86"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.05%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "235.06"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.22%  "

$ws.Range("E24").Value = "  -2.05%  "

$ws.Range("E25").Value = "  +1.15%  "

$ws.Range("E26").Value = "  -0.02%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "24.33"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.54%  "

$ws.Range("E28").Value = "  +14.89%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "9.30"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +2.15%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "32.42"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +3.65%  "

$ws.Range("E31").Value = "  -0.04%  "

$ws.Range("E32").Value = "  +0.45%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "17.49"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.18%  "

$ws.Range("E34").Value = "  +4.30%  "

$ws.Range("E35").Value = "  +6.43%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "128.97"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -7.71%  "

$ws.Range("E37").Value = "  +0.45%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "4.32"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.45%  "

$ws.Range("E39").Value = "  +3.94%  "

$ws.Range("E40").Value = "  -2.76%  "

$ws.Range("E41").Value = "  -0.55%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "20.96"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -6.12%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.926.16"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.33%  "

$ws.Range("E44").Value = "  +0.00%  "

$ws.Range("E45").Value = "  +3.31%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "9.27"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -8.98%  "

$ws.Range("E47").Value = "  +0.50%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "2.586.79"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +1.37%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.50"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +2.61%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "71.36"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -1.08%  "

$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "51.10"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -2.95%  "
